# Adds the "I0" and "IF" columns (I and J) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
# Match the header formatting used by the other header cells (B1:H1) —
# copy H1's formatting (bold, centered, thin border) onto I1:J1 without
# disturbing the values we just wrote.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (2-37) ---
$data = @(
    @{Row=2; I=8; J=8}
    @{Row=3; I=8; J=8}
    @{Row=4; I=6; J=7}
    @{Row=5; I=9; J=9}
    @{Row=6; I=8; J=8}
    @{Row=7; I=7; J=7}
    @{Row=8; I=9; J=9}
    @{Row=9; I=7; J=7}
    @{Row=10; I=8; J=8}
    @{Row=11; I=7; J=7}
    @{Row=12; I=8; J=8}
    @{Row=13; I=7; J=7}
    @{Row=14; I=8; J=8}
    @{Row=15; I=6; J=7}
    @{Row=16; I=7; J=8}
    @{Row=17; I=9; J=9}
    @{Row=18; I=9; J=9}
    @{Row=19; I=8; J=8}
    @{Row=20; I=8; J=8}
    @{Row=21; I=3; J=4}
    @{Row=22; I=7; J=8}
    @{Row=23; I=4; J=4}
    @{Row=24; I=9; J=9}
    @{Row=25; I=9; J=9}
    @{Row=26; I=7; J=7}
    @{Row=27; I=7; J=8}
    @{Row=28; I=5; J=6}
    @{Row=29; I=7; J=8}
    @{Row=30; I=7; J=7}
    @{Row=31; I=7; J=8}
    @{Row=32; I=6; J=7}
    @{Row=33; I=5; J=5}
    @{Row=34; I=2; J=3}
    @{Row=35; I=6; J=6}
    @{Row=36; I=6; J=6}
    @{Row=37; I=4; J=4}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 9).Value = $entry.I   # column I
    $ws.Cells.Item($r, 10).Value = $entry.J  # column J
}
